$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.074.84"
$ws.Range("E2").Value = "  -7.32%  "

$ws.Range("D3").Value = "1.415.94"
$ws.Range("E3").Value = "  -7.82%  "

$ws.Range("D4").Value = "'0.9931"
$ws.Range("E4").Value = "  -0.60%  "

$ws.Range("D5").Value = "'0.9938"
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("D6").Value = "'272.56"
$ws.Range("E6").Value = "  -5.98%  "

$ws.Range("D7").Value = "'0.3684"
$ws.Range("E7").Value = "  -6.45%  "

$ws.Range("D8").Value = "'0.3141"
$ws.Range("E8").Value = "  -1.08%  "

$ws.Range("D9").Value = "'39.63"
$ws.Range("E9").Value = "  -6.51%  "

$ws.Range("D10").Value = "'1.020"
$ws.Range("E10").Value = "  -4.48%  "

$ws.Range("D11").Value = "'0.06516"
$ws.Range("E11").Value = "  -9.28%  "

$ws.Range("D12").Value = "'0.9925"
$ws.Range("E12").Value = "  -0.67%  "

$ws.Range("D13").Value = "'5.478"
$ws.Range("E13").Value = "  -4.16%  "

$ws.Range("D14").Value = "'17.69"
$ws.Range("E14").Value = "  -3.85%  "

$ws.Range("D15").Value = "'6.179"
$ws.Range("E15").Value = "  -6.82%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.00001022"
$ws.Range("E16").Value = "  -6.75%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.407.84"
$ws.Range("E17").Value = "  -8.32%  "

$ws.Range("D18").Value = "'0.05698"
$ws.Range("E18").Value = "  -13.74%  "

$ws.Range("D19").Value = "'71.79"
$ws.Range("E19").Value = "  -14.76%  "

$ws.Range("D20").Value = "'0.9937"
$ws.Range("E20").Value = "  -0.57%  "

$ws.Range("D21").Value = "'5.601"
$ws.Range("E21").Value = "  -9.00%  "

$ws.Range("D22").Value = "'14.95"
$ws.Range("E22").Value = "  -3.93%  "

$ws.Range("D23").Value = "'11.03"
$ws.Range("E23").Value = "  +2.65%  "

$ws.Range("D24").Value = "'2.257"
$ws.Range("E24").Value = "  -5.05%  "

$ws.Range("D25").Value = "20.099.89"
$ws.Range("E25").Value = "  -7.21%  "

$ws.Range("D26").Value = "'2.279"
$ws.Range("E26").Value = "  -3.99%  "

$ws.Range("D27").Value = "'136.46"
$ws.Range("E27").Value = "  -9.57%  "

$ws.Range("D28").Value = "'17.02"
$ws.Range("E28").Value = "  -7.35%  "

$ws.Range("D29").Value = "1.565.39"
$ws.Range("E29").Value = "  -8.24%  "

$ws.Range("D30").Value = "'109.97"
$ws.Range("E30").Value = "  -6.42%  "

$ws.Range("D31").Value = "'4.117"
$ws.Range("E31").Value = "  -15.15%  "

$ws.Range("D32").Value = "'5.346"
$ws.Range("E32").Value = "  -12.31%  "

$ws.Range("D33").Value = "'0.8409"
$ws.Range("E33").Value = "  -11.99%  "

$ws.Range("D34").Value = "'0.07668"
$ws.Range("E34").Value = "  -5.94%  "

$ws.Range("D35").Value = "'8.370"
$ws.Range("E35").Value = "  -1.64%  "

$ws.Range("D36").Value = "'1.458"
$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("D37").Value = "'0.05797"
$ws.Range("E37").Value = "  -3.62%  "

$ws.Range("D38").Value = "'4.836"
$ws.Range("E38").Value = "  -7.19%  "

$ws.Range("D39").Value = "'0.9941"
$ws.Range("E39").Value = "  -0.52%  "

$ws.Range("D40").Value = "'0.02080"
$ws.Range("E40").Value = "  -6.64%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'10.59"
$ws.Range("E41").Value = "  -4.80%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1913"
$ws.Range("E42").Value = "  -6.38%  "

$ws.Range("D43").Value = "'1.084"
$ws.Range("E43").Value = "  -8.47%  "

$ws.Range("D44").Value = "'0.5322"
$ws.Range("E44").Value = "  -8.78%  "

$ws.Range("D45").Value = "'12.34"
$ws.Range("E45").Value = "  -6.18%  "

$ws.Range("D46").Value = "'3.505"
$ws.Range("E46").Value = "  -5.90%  "

$ws.Range("D47").Value = "'0.5180"
$ws.Range("E47").Value = "  -7.00%  "

$ws.Range("D48").Value = "'112.16"
$ws.Range("E48").Value = "  -3.80%  "

$ws.Range("D49").Value = "'1.781"
$ws.Range("E49").Value = "  -5.64%  "

$ws.Range("D50").Value = "'1.042"
$ws.Range("E50").Value = "  -11.10%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.06194"
$ws.Range("E51").Value = "  -7.70%  "
